$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2578.8
$ws.Range("I18").Value = 2578.8
$ws.Range("K18").Value = 2578.8
$ws.Range("M18").Value = -2294.8
$ws.Range("H19").Value = 972
$ws.Range("I19").Value = 973
$ws.Range("J19").Value = 971
$ws.Range("K19").Value = 973
$ws.Range("L19").Value = 971
$ws.Range("M19").Value = -798
$ws.Range("N19").Value = -1321
$ws.Range("H32").Value = 7002
$ws.Range("J32").Value = 7002
$ws.Range("L32").Value = 7002
$ws.Range("N32").Value = -7654
$ws.Range("H70").Value = 11080.363
$ws.Range("I70").Value = 8748
$ws.Range("K70").Value = 26244
$ws.Range("M70").Value = -25974
$ws.Range("H73").Value = 11080.363
$ws.Range("I73").Value = 8748
$ws.Range("K73").Value = 26244
$ws.Range("M73").Value = -25308
$ws.Range("H74").Value = 150435
$ws.Range("I74").Value = 200978.28
$ws.Range("K74").Value = 200978.28
$ws.Range("M74").Value = -200042.28
$ws.Range("H77").Value = 150435
$ws.Range("I77").Value = 200978.28
$ws.Range("K77").Value = 1004891.4
$ws.Range("M77").Value = -1000211.4
$ws.Range("H111").Value = 2327.75
$ws.Range("I111").Value = 375.33334
$ws.Range("K111").Value = 1126.00002
$ws.Range("M111").Value = 1940.99998
$ws.Range("H113").Value = 3201.8572
$ws.Range("I113").Value = 2788.2
$ws.Range("K113").Value = 2788.2
$ws.Range("M113").Value = 465.8000000000002
$ws.Range("H125").Value = 750000000
$ws.Range("I125").Value = 1000000000
$ws.Range("K125").Value = 9000000000
$ws.Range("M125").Value = -8999997540
$ws.Range("H127").Value = 998
$ws.Range("I127").Value = 998
$ws.Range("K127").Value = 2994
$ws.Range("M127").Value = 1966

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 567.7273
$ws.Range("I2").Value = 374.5
$ws.Range("K2").Value = 374.5
$ws.Range("M2").Value = -261.5
$ws.Range("H116").Value = 567.7273
$ws.Range("I116").Value = 374.5
$ws.Range("K116").Value = 374.5
$ws.Range("M116").Value = 1919.5
$ws.Range("H132").Value = 1771.4667
$ws.Range("I132").Value = 1771.4667
$ws.Range("K132").Value = 5314.4001
$ws.Range("M132").Value = -2784.4001

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 567.7273
$ws.Range("I3").Value = 374.5
$ws.Range("K3").Value = 374.5
$ws.Range("M3").Value = -260.5
$ws.Range("H22").Value = 635.44446
$ws.Range("J22").Value = 969.5
$ws.Range("L22").Value = 969.5
$ws.Range("N22").Value = -1315.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 15166.883
$ws.Range("I22").Value = 509.5
$ws.Range("K22").Value = 509.5
$ws.Range("M22").Value = -159.5
$ws.Range("H29").Value = 19990
$ws.Range("J29").Value = 19990
$ws.Range("L29").Value = 19990
$ws.Range("N29").Value = -20576
$ws.Range("H31").Value = 1684.1578
$ws.Range("I31").Value = 1417.2667
$ws.Range("K31").Value = 1417.2667
$ws.Range("M31").Value = -1122.2667
$ws.Range("H34").Value = 1684.1578
$ws.Range("I34").Value = 1417.2667
$ws.Range("K34").Value = 1417.2667
$ws.Range("M34").Value = -1215.2667
$ws.Range("H99").Value = 3621.3333
$ws.Range("I99").Value = 1900
$ws.Range("J99").Value = 3965.6
$ws.Range("K99").Value = 1900
$ws.Range("L99").Value = 3965.6
$ws.Range("M99").Value = -402
$ws.Range("N99").Value = -6961.6
$ws.Range("H126").Value = 3621.3333
$ws.Range("I126").Value = 1900
$ws.Range("J126").Value = 3965.6
$ws.Range("K126").Value = 5700
$ws.Range("L126").Value = 11896.8
$ws.Range("M126").Value = -3230
$ws.Range("N126").Value = -16836.8
$ws.Range("H133").Value = 45219
$ws.Range("J133").Value = 46523.75
$ws.Range("L133").Value = 46523.75
$ws.Range("N133").Value = -51583.75
$ws.Range("H138").Value = 168281.5
$ws.Range("I138").Value = 20774
$ws.Range("J138").Value = 315789
$ws.Range("K138").Value = 20774
$ws.Range("L138").Value = 315789
$ws.Range("M138").Value = -15634
$ws.Range("N138").Value = -326069
$ws.Range("H141").Value = 1000000
$ws.Range("J141").Value = 1000000
$ws.Range("L141").Value = 1000000
$ws.Range("N141").Value = -1010360

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 158090
$ws.Range("I2").Value = 137617.5
$ws.Range("J2").Value = 185386.67
$ws.Range("K2").Value = 825705
$ws.Range("L2").Value = 1112320.02
$ws.Range("M2").Value = -825592
$ws.Range("N2").Value = -1112546.02
$ws.Range("H80").Value = 2180
$ws.Range("I80").Value = 2180
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 6540
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -5604
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 2180
$ws.Range("I83").Value = 2180
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 19620
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -14940
$ws.Range("N83").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H70").Value = 8662.77
$ws.Range("J70").Value = 8546.143
$ws.Range("L70").Value = 8546.143
$ws.Range("N70").Value = -9086.143
$ws.Range("H73").Value = 8662.77
$ws.Range("J73").Value = 8546.143
$ws.Range("L73").Value = 8546.143
$ws.Range("N73").Value = -10418.143
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H132").Value = 992
$ws.Range("I132").Value = 992.6667
$ws.Range("K132").Value = 2978.0001
$ws.Range("M132").Value = -448.0001000000002

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1188.6923
$ws.Range("I16").Value = 909.7273
$ws.Range("J16").Value = 2723
$ws.Range("K16").Value = 909.7273
$ws.Range("L16").Value = 2723
$ws.Range("M16").Value = -739.7273
$ws.Range("N16").Value = -3063
$ws.Range("H55").Value = 362.92856
$ws.Range("I55").Value = 314.07693
$ws.Range("K55").Value = 314.07693
$ws.Range("M55").Value = -141.07693
$ws.Range("H82").Value = 1759.7142
$ws.Range("I82").Value = 1783.6
$ws.Range("J82").Value = 1700
$ws.Range("K82").Value = 1783.6
$ws.Range("L82").Value = 1700
$ws.Range("M82").Value = -1422.6
$ws.Range("N82").Value = -2422
$ws.Range("H85").Value = 1759.7142
$ws.Range("I85").Value = 1783.6
$ws.Range("J85").Value = 1700
$ws.Range("K85").Value = 1783.6
$ws.Range("L85").Value = 1700
$ws.Range("M85").Value = -535.5999999999999
$ws.Range("N85").Value = -4196
$ws.Range("H122").Value = 8027.364
$ws.Range("I122").Value = 8891.833000000001
$ws.Range("J122").Value = 6990
$ws.Range("K122").Value = 26675.499
$ws.Range("L122").Value = 20970
$ws.Range("M122").Value = -24225.499
$ws.Range("N122").Value = -25870
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9949.75
$ws.Range("J41").Value = 9299.666999999999
$ws.Range("L41").Value = 9299.666999999999
$ws.Range("N41").Value = -10079.667
$ws.Range("H62").Value = 14995.8
$ws.Range("I62").Value = 14995
$ws.Range("J62").Value = 14996
$ws.Range("K62").Value = 14995
$ws.Range("L62").Value = 14996
$ws.Range("M62").Value = -14371
$ws.Range("N62").Value = -16244
$ws.Range("H65").Value = 14995.8
$ws.Range("I65").Value = 14995
$ws.Range("J65").Value = 14996
$ws.Range("K65").Value = 74975
$ws.Range("L65").Value = 74980
$ws.Range("M65").Value = -71855
$ws.Range("N65").Value = -81220
$ws.Range("H81").Value = 2001088
$ws.Range("I81").Value = 1359.5
$ws.Range("K81").Value = 2719
$ws.Range("M81").Value = -1658
$ws.Range("H84").Value = 2001088
$ws.Range("I84").Value = 1359.5
$ws.Range("K84").Value = 13595
$ws.Range("M84").Value = -8291
$ws.Range("H122").Value = 278.4
$ws.Range("I122").Value = 260.5
$ws.Range("J122").Value = 350
$ws.Range("K122").Value = 781.5
$ws.Range("L122").Value = 1050
$ws.Range("M122").Value = 1668.5
$ws.Range("N122").Value = -5950
$ws.Range("H126").Value = 4318.7144
$ws.Range("I126").Value = 2887.3125
$ws.Range("K126").Value = 8661.9375
$ws.Range("M126").Value = -6191.9375
